$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 71428744
$ws.Range("I11").Value = 71428744
$ws.Range("K11").Value = 71428744
$ws.Range("M11").Value = -71428604
$ws.Range("H32").Value = 5855.6665
$ws.Range("I32").Value = 6900
$ws.Range("J32").Value = 5646.8
$ws.Range("K32").Value = 6900
$ws.Range("L32").Value = 5646.8
$ws.Range("M32").Value = -6574
$ws.Range("N32").Value = -6298.8
$ws.Range("H46").Value = 19999.5
$ws.Range("J46").Value = 19999.5
$ws.Range("L46").Value = 59998.5
$ws.Range("N46").Value = -60236.5
$ws.Range("H60").Value = 19999.5
$ws.Range("J60").Value = 19999.5
$ws.Range("L60").Value = 59998.5
$ws.Range("N60").Value = -60966.5
$ws.Range("H80").Value = 3736941.2
$ws.Range("I80").Value = 1635528.8
$ws.Range("J80").Value = 7939766
$ws.Range("K80").Value = 4906586.4
$ws.Range("L80").Value = 23819298
$ws.Range("M80").Value = -4905588.4
$ws.Range("N80").Value = -23821294
$ws.Range("H83").Value = 3736941.2
$ws.Range("I83").Value = 1635528.8
$ws.Range("J83").Value = 7939766
$ws.Range("K83").Value = 14719759.2
$ws.Range("L83").Value = 71457894
$ws.Range("M83").Value = -14714767.2
$ws.Range("N83").Value = -71467878
$ws.Range("H92").Value = 4157.6
$ws.Range("I92").Value = 1426.1111
$ws.Range("K92").Value = 1426.1111
$ws.Range("M92").Value = -178.1111000000001
$ws.Range("H137").Value = 2356.5217
$ws.Range("I137").Value = 1129.2858
$ws.Range("J137").Value = 4265.5557
$ws.Range("K137").Value = 3387.8574
$ws.Range("L137").Value = 12796.6671
$ws.Range("M137").Value = -837.8574000000003
$ws.Range("N137").Value = -17896.6671

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 881.61536
$ws.Range("J5").Value = 1631.3334
$ws.Range("L5").Value = 1631.3334
$ws.Range("N5").Value = -1855.3334
$ws.Range("H55").Value = 46667.332
$ws.Range("I55").Value = 39999.5
$ws.Range("J55").Value = 60003
$ws.Range("K55").Value = 39999.5
$ws.Range("L55").Value = 60003
$ws.Range("M55").Value = -39684.5
$ws.Range("N55").Value = -60633
$ws.Range("H60").Value = 79203.375
$ws.Range("I60").Value = 79203.375
$ws.Range("K60").Value = 79203.375
$ws.Range("M60").Value = -78470.375

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 881.61536
$ws.Range("J4").Value = 1631.3334
$ws.Range("L4").Value = 1631.3334
$ws.Range("N4").Value = -1861.3334
$ws.Range("H107").Value = 2643.2222
$ws.Range("I107").Value = 2348.625
$ws.Range("K107").Value = 2348.625
$ws.Range("M107").Value = -428.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27029330
$ws.Range("J31").Value = 4538.2
$ws.Range("L31").Value = 4538.2
$ws.Range("N31").Value = -5128.2
$ws.Range("H34").Value = 27029330
$ws.Range("J34").Value = 4538.2
$ws.Range("L34").Value = 4538.2
$ws.Range("N34").Value = -4942.2
$ws.Range("H58").Value = 3079.963
$ws.Range("J58").Value = 3200.5386
$ws.Range("L58").Value = 3200.5386
$ws.Range("N58").Value = -3606.5386
$ws.Range("H105").Value = 1434.7333
$ws.Range("I105").Value = 1152
$ws.Range("J105").Value = 2000.2
$ws.Range("K105").Value = 1152
$ws.Range("L105").Value = 2000.2
$ws.Range("M105").Value = 595
$ws.Range("N105").Value = -5494.2
$ws.Range("H107").Value = 1087.1818
$ws.Range("I107").Value = 477
$ws.Range("J107").Value = 1819.4
$ws.Range("K107").Value = 477
$ws.Range("L107").Value = 1819.4
$ws.Range("M107").Value = 1443
$ws.Range("N107").Value = -5659.4
$ws.Range("H111").Value = 83566.664
$ws.Range("J111").Value = 83566.664
$ws.Range("L111").Value = 83566.664
$ws.Range("N111").Value = -91746.664
$ws.Range("H136").Value = 3079.963
$ws.Range("J136").Value = 3200.5386
$ws.Range("L136").Value = 9601.6158
$ws.Range("N136").Value = -14701.6158

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 774.5
$ws.Range("I5").Value = 774.5
$ws.Range("K5").Value = 2323.5
$ws.Range("M5").Value = -2211.5
$ws.Range("H34").Value = 7080.2144
$ws.Range("J34").Value = 7080.2144
$ws.Range("L34").Value = 21240.6432
$ws.Range("N34").Value = -21408.6432
$ws.Range("H76").Value = 16674
$ws.Range("I76").Value = 15
$ws.Range("K76").Value = 45
$ws.Range("M76").Value = 338
$ws.Range("H79").Value = 16674
$ws.Range("I79").Value = 15
$ws.Range("K79").Value = 45
$ws.Range("M79").Value = 1281
$ws.Range("H114").Value = 5353.4443
$ws.Range("J114").Value = 5353.4443
$ws.Range("L114").Value = 16060.3329
$ws.Range("N114").Value = -22568.3329
$ws.Range("H135").Value = 774.5
$ws.Range("I135").Value = 774.5
$ws.Range("K135").Value = 6970.5
$ws.Range("M135").Value = -4435.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 8757.6
$ws.Range("I97").Value = 1643.8
$ws.Range("J97").Value = 12314.5
$ws.Range("K97").Value = 1643.8
$ws.Range("L97").Value = 12314.5
$ws.Range("M97").Value = -1147.8
$ws.Range("N97").Value = -13306.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7186.4287
$ws.Range("I40").Value = 7186.4287
$ws.Range("K40").Value = 7186.4287
$ws.Range("M40").Value = -7050.4287
$ws.Range("H46").Value = 1435.5
$ws.Range("I46").Value = 498
$ws.Range("K46").Value = 498
$ws.Range("M46").Value = -310
$ws.Range("H57").Value = 33826.332
$ws.Range("I57").Value = 27810.637
$ws.Range("K57").Value = 27810.637
$ws.Range("M57").Value = -27244.637
$ws.Range("H68").Value = 2780282.5
$ws.Range("I68").Value = 4631688
$ws.Range("K68").Value = 4631688
$ws.Range("M68").Value = -4630939
$ws.Range("H71").Value = 2780282.5
$ws.Range("I71").Value = 4631688
$ws.Range("K71").Value = 23158440
$ws.Range("M71").Value = -23154696
$ws.Range("H100").Value = 20859944
$ws.Range("I100").Value = 4897.3335
$ws.Range("J100").Value = 41714990
$ws.Range("K100").Value = 4897.3335
$ws.Range("L100").Value = 41714990
$ws.Range("M100").Value = -4356.3335
$ws.Range("N100").Value = -41716072
$ws.Range("H122").Value = 4481.204
$ws.Range("J122").Value = 8688.299999999999
$ws.Range("L122").Value = 26064.9
$ws.Range("N122").Value = -30964.9
